$wb = $excel.ActiveWorkbook

# Work on the second sheet (TestData2), which gets the new data table and becomes the active/printed sheet
$ws2 = $wb.Worksheets.Item("TestData2")

$ws2.Range("A1").Value = "Value1"
$ws2.Range("B1").Value = "Value2"
$ws2.Range("C1").Value = "Value3"
$ws2.Range("D1").Value = "Value4"
$ws2.Range("E1").Value = "Value5"

$ws2.Range("A2").Value = "a"
$ws2.Range("A3").Value = "b"
$ws2.Range("A4").Value = "c"
$ws2.Range("A5").Value = "d"

$ws2.Range("B2").Value = "e"
$ws2.Range("B3").Value = "f"
$ws2.Range("B4").Value = "g"
$ws2.Range("B5").Value = "h"

$ws2.Range("C2").Value = "i"
$ws2.Range("C3").Value = "j"
$ws2.Range("C4").Value = "k"
$ws2.Range("C5").Value = "l"

$ws2.Range("D2").Value = "m"
$ws2.Range("D3").Value = "n"
$ws2.Range("D4").Value = "o"
$ws2.Range("D5").Value = "p"

$ws2.Range("E2").Value = "q"
$ws2.Range("E3").Value = "r"
$ws2.Range("E4").Value = "s"
$ws2.Range("E5").Value = "t"

# Select F9 on TestData2 to match the saved selection/active cell, then activate the sheet
$ws2.Range("F9").Select()
$ws2.Activate()

$wb.Save()
